# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-08 (serial 45177) to 2023-09-09 (serial 45178).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 9 -Day 9 -Hour 0 -Minute 0 -Second 0 -Millisecond 0

$range = $ws.Range("C2:C246")
$range.Value = $newDate
